$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hydro")

# Update source data values for the Hydro (UK) section
$ws.Range("C10").Value = 68317
$ws.Range("C13").Value = 3749.0329999999972
$ws.Range("C14").Value = 249

# Update selection to match the saved view state (D13 active cell)
$ws.Range("D13").Select()
